$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 996
$ws.Range("B2").Value = 981
$ws.Range("C2").Value = 981
$ws.Range("D2").Value = 981
$ws.Range("G2").Value = 996
